$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (simplified) Sample ID query text - replaces the old B3 content which had
# extra Tumor / Analyte Type columns. The old text is dropped entirely from the
# shared string table and this new text is appended as a new shared string.
$newSampleQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001819' AND gi.instrument_model = 'Illumina HiSeq 2000'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

# Update the SamplesTab query cell (B3) with the new, simplified query.
$ws.Range("B3").Value = $newSampleQuery

# Move the cursor/selection to match the saved view state (B3 selected,
# top-left cell scrolled to A3).
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 3
